# Updates the cryptos price/volume table (and a few re-ranked rows) to the
# latest scraped values. The "Price" column (D) is stored as plain text in
# this workbook (e.g. "42.930.45", "34.09") even when a value happens to
# look like a number, so purely-numeric replacements are written with a
# leading apostrophe - Excel's standard "force text" prefix - to stop them
# from being auto-coerced into numbers (which would silently drop trailing
# zeros, e.g. "126.60" -> 126.6, or introduce floating point artifacts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.921.50'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '2.361.11'
$ws.Range("E3").Value = '  +2.24%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''302.15'

$ws.Range("D6").Value = '''95.75'
$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").Value = '''34.08'
$ws.Range("E10").Value = '  -0.17%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.124'
$ws.Range("E11").Value = '  +3.52%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '''0.0783'
$ws.Range("E12").Value = '  +0.18%  '

$ws.Range("D13").Value = '''18.29'
$ws.Range("E13").Value = '  -2.93%  '

$ws.Range("D14").Value = '''6.71'
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").Value = '2.729.95'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").Value = '2.366.85'
$ws.Range("E16").Value = '  +3.81%  '

$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '42.890.42'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("E19").Value = '  +1.96%  '

$ws.Range("E20").Value = '  -2.65%  '

$ws.Range("E21").Value = '  -0.85%  '

$ws.Range("D22").Value = '''67.84'
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").Value = '''234.85'
$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("E24").Value = '  -4.73%  '

$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  +0.83%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").Value = '''9.28'
$ws.Range("E29").Value = '  +2.10%  '

$ws.Range("D30").Value = '''31.88'
$ws.Range("E30").Value = '  -0.61%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("D33").Value = '''17.27'
$ws.Range("E33").Value = '  -1.99%  '

$ws.Range("D34").Value = '''0.0709'
$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("E35").Value = '  +3.65%  '

$ws.Range("E36").Value = '  +3.44%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''126.60'
$ws.Range("E37").Value = '  -23.95%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''4.33'
$ws.Range("E38").Value = '  -2.85%  '

$ws.Range("D39").Value = '''2.28'
$ws.Range("E39").Value = '  -2.39%  '

$ws.Range("D40").Value = '''2.79'
$ws.Range("E40").Value = '  +3.64%  '

$ws.Range("E41").Value = '  -0.78%  '

$ws.Range("D42").Value = '''21.22'

$ws.Range("D43").Value = '1.927.83'
$ws.Range("E43").Value = '  +0.26%  '

$ws.Range("E44").Value = '  -0.23%  '

$ws.Range("E45").Value = '  +2.35%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''9.16'
$ws.Range("E46").Value = '  -8.64%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.71'
$ws.Range("E47").Value = '  -0.66%  '

$ws.Range("D48").Value = '2.594.95'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("E50").Value = '  +1.56%  '

$ws.Range("D51").Value = '''51.46'
$ws.Range("E51").Value = '  -3.41%  '
